$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Change "10 provinsi" -> "5 provinsi" in A3
$ws.Range("A3").Value = "5 provinsi"

# 2. Change selection to C2
$ws.Range("C2").Select()

# 3. Center-align B4:N8 (matches style index 1 instead of 2)
$ws.Range("B4:N8").HorizontalAlignment = -4108

# 4. Row 9 becomes a new "Total" row (was previously row 14)
$ws.Range("A9").Value = "Total "
$ws.Range("B9:M9").HorizontalAlignment = -4108
$ws.Range("N9").HorizontalAlignment = -4108
$ws.Range("B9").Formula = "=SUM(B4:B8)"
$ws.Range("C9:M9").Formula = "=SUM(C4:C8)"

# 5. Clear the old totals row (row 14): A14/B14 are fully removed (content + format),
#    C14:N14 keep their original style but lose their formula/value.
$ws.Range("A14:B14").Clear()
$ws.Range("C14:N14").ClearContents()

Write-Host "done"
